$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''51.579.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").Value = '''2.790.90'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''354.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = '''108.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("D9").Value = '''0.623'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.80%  '
$ws.Range("D10").Value = '''39.93'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("E11").Value = '  +0.92%  '
$ws.Range("D12").Value = '''0.0838'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("D13").Value = '''20.02'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.68%  '
$ws.Range("D14").Value = '''7.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.17%  '
$ws.Range("D15").Value = '''3.233.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("D16").Value = '''2.788.95'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '''51.545.39'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("E19").Value = '  +3.78%  '
$ws.Range("D20").Value = '''3.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.71%  '
$ws.Range("E21").Value = '  +1.97%  '
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("D23").Value = '''70.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '''267.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.19%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  -1.60%  '
$ws.Range("D28").Value = '''0.166'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").Value = '''10.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("D30").Value = '''37.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.81%  '
$ws.Range("E31").Value = '  +4.40%  '
$ws.Range("D32").Value = '''6.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.03%  '
$ws.Range("D33").Value = '''51.87'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("E34").Value = '  +9.43%  '
$ws.Range("E35").Value = '  -6.17%  '
$ws.Range("E36").Value = '  +0.82%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '''18.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("E39").Value = '  -2.13%  '
$ws.Range("D40").Value = '''1.98'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").Value = '''2.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.93%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '''119.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '''21.76'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.87%  '
$ws.Range("D45").Value = '''2.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.88%  '
$ws.Range("D46").Value = '''2.126.24'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.21%  '
$ws.Range("D47").Value = '''3.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.85%  '
$ws.Range("E48").Value = '  +6.49%  '
$ws.Range("E49").Value = '  +10.72%  '
$ws.Range("E50").Value = '  -4.18%  '
$ws.Range("D51").Value = '''5.37'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.67%  '
